$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.192.22"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.601.86"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.10"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3781"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.86"
$ws.Range("E8").Value = "  +3.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3613"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.263"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08125"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.56"
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.582"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.383"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001247"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.603.14"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.72"
$ws.Range("E18").Value = "  +2.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06864"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.01"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.532"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.198.02"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.393"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.983"
$ws.Range("E26").Value = "  +8.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.17"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.237"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.63"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.823"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.779.18"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9769"
$ws.Range("E34").Value = "  +2.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07524"
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02723"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.125"
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2501"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08792"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7096"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.360"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.44"
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.44"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6538"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.304"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.015"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.26"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07953"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.203"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.229"
$ws.Range("E51").Value = "  +3.58%  "
